$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 340 (existing rows 340.. shift down to 342..)
$ws.Range("A340:A341").EntireRow.Insert()

# ---- New row 340 ----
$ws.Range("A340").Value = 6
$ws.Range("B340").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C340").Value = "Metropolitana"
$ws.Range("D340").Value = 45275
$ws.Range("E340").Value = 13
$ws.Range("F340").Value = "Fruta"
$ws.Range("G340").Value = 100101
$ws.Range("H340").Value = "Berries"
$ws.Range("I340").Value = 100101004
$ws.Range("J340").Value = "Frambuesa"
$ws.Range("K340").Value = "Sin especificar"
$ws.Range("L340").Value = "Especial"
$ws.Range("M340").Value = 250
$ws.Range("N340").Value = 9000
$ws.Range("O340").Value = 9000
$ws.Range("P340").Value = 9000
$ws.Range("Q340").Value = "$/bandeja 2 kilos"
$ws.Range("R340").Value = "Provincia de Curicó"
$ws.Range("S340").Value = 4500
$ws.Range("T340").Value = 2

# ---- New row 341 ----
$ws.Range("A341").Value = 6
$ws.Range("B341").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C341").Value = "Metropolitana"
$ws.Range("D341").Value = 45275
$ws.Range("E341").Value = 13
$ws.Range("F341").Value = "Fruta"
$ws.Range("G341").Value = 100101
$ws.Range("H341").Value = "Berries"
$ws.Range("I341").Value = 100101004
$ws.Range("J341").Value = "Frambuesa"
$ws.Range("K341").Value = "Sin especificar"
$ws.Range("L341").Value = "Primera"
$ws.Range("M341").Value = 300
$ws.Range("N341").Value = 8000
$ws.Range("O341").Value = 8000
$ws.Range("P341").Value = 8000
$ws.Range("Q341").Value = "$/bandeja 2 kilos"
$ws.Range("R341").Value = "Provincia de Curicó"
$ws.Range("S341").Value = 4000
$ws.Range("T341").Value = 2
